# Increment the "want to go" counts (column F) on the "展览" and
# "全部类型" sheets to reflect the newly generated totals:
#   F4: 96  -> 97
#   F5: 863 -> 866
#   F6: 209 -> 210

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 97
    $ws.Range("F5").Value = 866
    $ws.Range("F6").Value = 210
}
